# Auto-generated edit script: Natmi following Dr Hou advice
# Adds a new "ECs" sending-cluster block (2 rows) ahead of the existing
# FAPs/M2/sCs blocks, and recomputes all Cfh-Itgam LR-pair metrics so the
# sheet grows from 6 data rows (A1:T7) to 8 data rows (A1:T9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> M2
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Cfh"
$ws.Cells.Item(2, 3).Value2 = "Itgam"
$ws.Cells.Item(2, 4).Value2 = "M2"
$ws.Cells.Item(2, 5).Value2 = 2
$ws.Cells.Item(2, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(2, 7).Value2 = 0.6849806666666667
$ws.Cells.Item(2, 8).Value2 = 2.054942
$ws.Cells.Item(2, 9).Value2 = 0.005384415753505337
$ws.Cells.Item(2, 10).Value2 = 0.005384415753505337
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 45.931316
$ws.Cells.Item(2, 14).Value2 = 137.793948
$ws.Cells.Item(2, 15).Value2 = 0.9874217014725413
$ws.Cells.Item(2, 16).Value2 = 0.9874217014725412
$ws.Cells.Item(2, 17).Value2 = 31.46206345455733
$ws.Cells.Item(2, 18).Value2 = 283.158571091016
$ws.Cells.Item(2, 19).Value2 = 0.005316688964761796
$ws.Cells.Item(2, 20).Value2 = 0.005316688964761795

# Row 3: ECs -> sCs
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Cfh"
$ws.Cells.Item(3, 3).Value2 = "Itgam"
$ws.Cells.Item(3, 4).Value2 = "sCs"
$ws.Cells.Item(3, 5).Value2 = 2
$ws.Cells.Item(3, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(3, 7).Value2 = 0.6849806666666667
$ws.Cells.Item(3, 8).Value2 = 2.054942
$ws.Cells.Item(3, 9).Value2 = 0.005384415753505337
$ws.Cells.Item(3, 10).Value2 = 0.005384415753505337
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 0.5850973333333334
$ws.Cells.Item(3, 14).Value2 = 1.755292
$ws.Cells.Item(3, 15).Value2 = 0.01257829852745884
$ws.Cells.Item(3, 16).Value2 = 0.01257829852745884
$ws.Cells.Item(3, 17).Value2 = 0.4007803614515556
$ws.Cells.Item(3, 18).Value2 = 3.607023253064
$ws.Cells.Item(3, 19).Value2 = 0.00006772678874354236
$ws.Cells.Item(3, 20).Value2 = 0.00006772678874354236

# Row 4: FAPs -> M2
$ws.Cells.Item(4, 1).Value2 = "FAPs"
$ws.Cells.Item(4, 2).Value2 = "Cfh"
$ws.Cells.Item(4, 3).Value2 = "Itgam"
$ws.Cells.Item(4, 4).Value2 = "M2"
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 97.08941133333333
$ws.Cells.Item(4, 8).Value2 = 291.268234
$ws.Cells.Item(4, 9).Value2 = 0.7631890669640694
$ws.Cells.Item(4, 10).Value2 = 0.7631890669640695
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 45.931316
$ws.Cells.Item(4, 14).Value2 = 137.793948
$ws.Cells.Item(4, 15).Value2 = 0.9874217014725413
$ws.Cells.Item(4, 16).Value2 = 0.9874217014725412
$ws.Cells.Item(4, 17).Value2 = 4459.444432205315
$ws.Cells.Item(4, 18).Value2 = 40134.99988984783
$ws.Cells.Item(4, 19).Value2 = 0.7535894470469027
$ws.Cells.Item(4, 20).Value2 = 0.7535894470469027

# Row 5: FAPs -> sCs
$ws.Cells.Item(5, 1).Value2 = "FAPs"
$ws.Cells.Item(5, 2).Value2 = "Cfh"
$ws.Cells.Item(5, 3).Value2 = "Itgam"
$ws.Cells.Item(5, 4).Value2 = "sCs"
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 97.08941133333333
$ws.Cells.Item(5, 8).Value2 = 291.268234
$ws.Cells.Item(5, 9).Value2 = 0.7631890669640694
$ws.Cells.Item(5, 10).Value2 = 0.7631890669640695
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 0.5850973333333334
$ws.Cells.Item(5, 14).Value2 = 1.755292
$ws.Cells.Item(5, 15).Value2 = 0.01257829852745884
$ws.Cells.Item(5, 16).Value2 = 0.01257829852745884
$ws.Cells.Item(5, 17).Value2 = 56.80675566603644
$ws.Cells.Item(5, 18).Value2 = 511.260800994328
$ws.Cells.Item(5, 19).Value2 = 0.00959961991716684
$ws.Cells.Item(5, 20).Value2 = 0.009599619917166841

# Row 6: M2 -> M2
$ws.Cells.Item(6, 1).Value2 = "M2"
$ws.Cells.Item(6, 2).Value2 = "Cfh"
$ws.Cells.Item(6, 3).Value2 = "Itgam"
$ws.Cells.Item(6, 4).Value2 = "M2"
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 19.002885
$ws.Cells.Item(6, 8).Value2 = 57.008655
$ws.Cells.Item(6, 9).Value2 = 0.1493756515114056
$ws.Cells.Item(6, 10).Value2 = 0.1493756515114056
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 45.931316
$ws.Cells.Item(6, 14).Value2 = 137.793948
$ws.Cells.Item(6, 15).Value2 = 0.9874217014725413
$ws.Cells.Item(6, 16).Value2 = 0.9874217014725412
$ws.Cells.Item(6, 17).Value2 = 872.82751584666
$ws.Cells.Item(6, 18).Value2 = 7855.44764261994
$ws.Cells.Item(6, 19).Value2 = 0.1474967599739615
$ws.Cells.Item(6, 20).Value2 = 0.1474967599739615

# Row 7: M2 -> sCs
$ws.Cells.Item(7, 1).Value2 = "M2"
$ws.Cells.Item(7, 2).Value2 = "Cfh"
$ws.Cells.Item(7, 3).Value2 = "Itgam"
$ws.Cells.Item(7, 4).Value2 = "sCs"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 19.002885
$ws.Cells.Item(7, 8).Value2 = 57.008655
$ws.Cells.Item(7, 9).Value2 = 0.1493756515114056
$ws.Cells.Item(7, 10).Value2 = 0.1493756515114056
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 0.5850973333333334
$ws.Cells.Item(7, 14).Value2 = 1.755292
$ws.Cells.Item(7, 15).Value2 = 0.01257829852745884
$ws.Cells.Item(7, 16).Value2 = 0.01257829852745884
$ws.Cells.Item(7, 17).Value2 = 11.11853733914
$ws.Cells.Item(7, 18).Value2 = 100.06683605226
$ws.Cells.Item(7, 19).Value2 = 0.001878891537444118
$ws.Cells.Item(7, 20).Value2 = 0.001878891537444118

# Row 8: sCs -> M2
$ws.Cells.Item(8, 1).Value2 = "sCs"
$ws.Cells.Item(8, 2).Value2 = "Cfh"
$ws.Cells.Item(8, 3).Value2 = "Itgam"
$ws.Cells.Item(8, 4).Value2 = "M2"
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 10.43813466666667
$ws.Cells.Item(8, 8).Value2 = 31.314404
$ws.Cells.Item(8, 9).Value2 = 0.0820508657710196
$ws.Cells.Item(8, 10).Value2 = 0.0820508657710196
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 12).Value2 = 1
$ws.Cells.Item(8, 13).Value2 = 45.931316
$ws.Cells.Item(8, 14).Value2 = 137.793948
$ws.Cells.Item(8, 15).Value2 = 0.9874217014725413
$ws.Cells.Item(8, 16).Value2 = 0.9874217014725412
$ws.Cells.Item(8, 17).Value2 = 479.4372618252214
$ws.Cells.Item(8, 18).Value2 = 4314.935356426992
$ws.Cells.Item(8, 19).Value2 = 0.08101880548691527
$ws.Cells.Item(8, 20).Value2 = 0.08101880548691526

# Row 9: sCs -> sCs
$ws.Cells.Item(9, 1).Value2 = "sCs"
$ws.Cells.Item(9, 2).Value2 = "Cfh"
$ws.Cells.Item(9, 3).Value2 = "Itgam"
$ws.Cells.Item(9, 4).Value2 = "sCs"
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 10.43813466666667
$ws.Cells.Item(9, 8).Value2 = 31.314404
$ws.Cells.Item(9, 9).Value2 = 0.0820508657710196
$ws.Cells.Item(9, 10).Value2 = 0.0820508657710196
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 0.5850973333333334
$ws.Cells.Item(9, 14).Value2 = 1.755292
$ws.Cells.Item(9, 15).Value2 = 0.01257829852745884
$ws.Cells.Item(9, 16).Value2 = 0.01257829852745884
$ws.Cells.Item(9, 17).Value2 = 6.107324758440889
$ws.Cells.Item(9, 18).Value2 = 54.965922825968
$ws.Cells.Item(9, 19).Value2 = 0.001032060284104339
$ws.Cells.Item(9, 20).Value2 = 0.001032060284104339

